$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "_Main" sheet: collapse the per-status-level Python/Java implementation
# columns (B/C) for rows that never used them, append a trailing newline to
# the Minio-storage use case description, and append a brand new use case
# (#300 - Face recognition) describing the JSON contract for face detection
# output together with the python runner command.
# ---------------------------------------------------------------------------
$main = $wb.Worksheets.Item("_Main")

$main.Range("C2").Clear()
$main.Range("B4:C4").Clear()
$main.Range("B5:C5").Clear()
$main.Range("B6:C6").Clear()
$main.Range("B7:C7").Clear()
$main.Range("B8:C8").Clear()
$main.Range("B9:C9").Clear()
$main.Range("B10:C10").Clear()
$main.Range("B11:C11").Clear()

# Re-affirm row 11's description text (adds the trailing line break that the
# authors introduced when the text got re-wrapped into the shared string
# table).
$main.Range("D11").Value = "Processed image must be saved to Minio at identifier, described on 100-x-image-input`n"

# New use case row - #300, Face recognition. (Values are populated in this
# specific order - details, then run command, then summary - so the shared
# string table ends up in the same append order the original authors hit.)
$main.Range("A12").Value = 300
$faceJson = "Input: jpeg image`noutput: array`n[`n        {`n            `"detection`": 0.92,`n            `"faceBox`": {`n                 `"p1`": { `"x`": 507, `"y`": 42 },`n                 `"p2`": { `"x`": 601, `"y`": 165 }`n              }`n        }`n ]`n"
$main.Range("E12").Value = $faceJson
$main.Range("E12").WrapText = $true
$main.Range("F12").Value = "python FacesImageProcessor.py --file=out/orban_putin.jpg"
$main.Range("D12").Value = "Face recognition: input image should be processed to output result in defined format"

# Row heights re-flowed slightly (autofit under the refreshed default font
# metrics) once the sheet content changed.
$main.Rows.Item(1).RowHeight = 10.5
$main.Rows.Item(2).RowHeight = 30
$main.Rows.Item(3).RowHeight = 40
$main.Rows.Item(4).RowHeight = 20
$main.Rows.Item(5).RowHeight = 40
$main.Rows.Item(6).RowHeight = 30
$main.Rows.Item(7).RowHeight = 30
$main.Rows.Item(8).RowHeight = 30
$main.Rows.Item(9).RowHeight = 30
$main.Rows.Item(10).RowHeight = 30
$main.Rows.Item(11).RowHeight = 30
$main.Rows.Item(12).RowHeight = 120

# ---------------------------------------------------------------------------
# "100-x-image-input" sheet: only cosmetic row-height reflow - the shared
# string indices referenced by row 14 shift automatically once the string
# table above is edited, so no cell content needs to be touched here.
# ---------------------------------------------------------------------------
$imageInput = $wb.Worksheets.Item("100-x-image-input")
$imageInput.Rows.Item(7).RowHeight = 10.5
$imageInput.Rows.Item(14).RowHeight = 60
$imageInput.Rows.Item(16).RowHeight = 10.5
$imageInput.Rows.Item(17).RowHeight = 230

# ---------------------------------------------------------------------------
# "_ListValues" sheet: drop the custom height on row 1 so it falls back to
# the sheet default again.
# ---------------------------------------------------------------------------
$listValues = $wb.Worksheets.Item("_ListValues")
$listValues.Rows.Item(1).AutoFit()

# ---------------------------------------------------------------------------
# The "_Main" tab becomes the active / selected tab (it used to be
# "100-x-image-input").
# ---------------------------------------------------------------------------
$main.Activate()
$main.Range("D13").Select()
